$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

$data = @{
    2  = @(-0.04123248769265224, 0.5678141048133444, 0.7453484225297263, 0.8633356372406541, 0.8836462639047493, 21)
    3  = @(0.7593282539754347, 0.8681867713281692, 1.984319268858228, 1.408658677202618, 1.217304992087394, 20)
    4  = @(0.4143225526684406, 1.117115571250938, 3.371275995140762, 1.836103481599216, 1.837761977278195, 19)
    5  = @(0.5506979343386952, 0.7414509954072014, 0.7679265988674776, 0.8763142124075574, 0.7014208651723969, 18)
    6  = @(0.6306395624298597, 0.6531076179729335, 0.6775008370835029, 0.8231043901495745, 0.5452354909515147, 17)
    7  = @(0.3250633412777073, 0.5387876674859324, 0.413150377612954, 0.642767747800832, 0.5726981303924567, 16)
    8  = @(0.3851030143285183, 0.4236031049613062, 0.2321424284816434, 0.4818116109867459, 0.2997107477636902, 15)
    9  = @(0.3424193635590594, 0.3780005890311999, 0.1990162949752141, 0.4461124241435269, 0.2999029874451486, 11)
    10 = @(0.3499679482706219, 0.3499679482706219, 0.18669282775681, 0.4320796544120193, 0.2775937959106322, 6)
    11 = @(0.4882709805783148, 0.4882709805783148, 0.2805492020480977, 0.5296689551484943, 0.2514179336479061, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
